$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price text stays as text (matches source inline-string cells)
$forceTextCells = @('D5','D10','D11','D16','D18','D25','D28','D40','D42','D43','D44','D47','D49','D51')
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the crypto price refresh
$ws.Range('D2').Value = '26.966.10'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '1.655.34'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '215.14'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +2.57%  '
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('D10').Value = '20.13'
$ws.Range('E10').Value = '  +4.58%  '
$ws.Range('D11').Value = '0.0878'
$ws.Range('E11').Value = '  +3.51%  '
$ws.Range('D12').Value = '1.889.02'
$ws.Range('E12').Value = '  +2.88%  '
$ws.Range('D13').Value = '1.656.65'
$ws.Range('E13').Value = '  +2.81%  '
$ws.Range('E14').Value = '  +2.06%  '
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('D16').Value = '65.25'
$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D17').Value = '26.980.28'
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').Value = '236.88'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('E22').Value = '  +3.71%  '
$ws.Range('E23').Value = '  +3.02%  '
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('D25').Value = '145.21'
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').Value = '15.86'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('E31').Value = '  +2.03%  '
$ws.Range('D32').Value = '1.554.21'
$ws.Range('E32').Value = '  +4.09%  '
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('E34').Value = '  +4.32%  '
$ws.Range('E35').Value = '  +10.26%  '
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('E37').Value = '  +3.60%  '
$ws.Range('E38').Value = '  +9.01%  '
$ws.Range('E39').Value = '  +2.63%  '
$ws.Range('D40').Value = '6.01'
$ws.Range('E40').Value = '  +3.88%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = '66.84'
$ws.Range('E42').Value = '  +9.44%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '0.970'
$ws.Range('E43').Value = '  +4.50%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '2.23'
$ws.Range('E44').Value = '  +2.22%  '
$ws.Range('D45').Value = '1.796.95'
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('D47').Value = '90.07'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('D49').Value = '0.0992'
$ws.Range('E49').Value = '  +3.23%  '
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('D51').Value = '7.67'
$ws.Range('E51').Value = '  +2.74%  '
